# Correccion a Diebold Mariano y revision de Cap1
# Updates the P_valores and Estadisticos_DM sheets with corrected values.

$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.977952015510797
$wsP.Range("D2").Value = 0.9808388988468704
$wsP.Range("E2").Value = 0.8938345121703961
$wsP.Range("F2").Value = 0.9113812386444586

$wsP.Range("B3").Value = 0.977952015510797
$wsP.Range("D3").Value = 0.9330881026500384
$wsP.Range("E3").Value = 0.9004553961666277
$wsP.Range("F3").Value = 0.8923606465832117

$wsP.Range("B4").Value = 0.9808388988468704
$wsP.Range("C4").Value = 0.9330881026500384
$wsP.Range("E4").Value = 0.7350034373220531
$wsP.Range("F4").Value = 0.8180493084429803

$wsP.Range("B5").Value = 0.8938345121703961
$wsP.Range("C5").Value = 0.9004553961666277
$wsP.Range("D5").Value = 0.7350034373220531
$wsP.Range("F5").Value = 0.9776962784895031

$wsP.Range("B6").Value = 0.9113812386444586
$wsP.Range("C6").Value = 0.8923606465832117
$wsP.Range("D6").Value = 0.8180493084429803
$wsP.Range("E6").Value = 0.9776962784895031

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 0.02813448171661221
$wsE.Range("D2").Value = -0.02444980859885563
$wsE.Range("E2").Value = 0.1359007548887671
$wsE.Range("F2").Value = 0.1133258852422692

$wsE.Range("B3").Value = -0.02813448171661221
$wsE.Range("D3").Value = -0.08548271405322977
$wsE.Range("E3").Value = 0.1273745628569715
$wsE.Range("F3").Value = 0.1378001744382444

$wsE.Range("B4").Value = 0.02444980859885563
$wsE.Range("C4").Value = 0.08548271405322977
$wsE.Range("E4").Value = 0.3452980656664401
$wsE.Range("F4").Value = 0.2344242761364861

$wsE.Range("B5").Value = -0.1359007548887671
$wsE.Range("C5").Value = -0.1273745628569715
$wsE.Range("D5").Value = -0.3452980656664401
$wsE.Range("F5").Value = -0.0284609105736368

$wsE.Range("B6").Value = -0.1133258852422692
$wsE.Range("C6").Value = -0.1378001744382444
$wsE.Range("D6").Value = -0.2344242761364861
$wsE.Range("E6").Value = 0.0284609105736368
